# Apply updated Plan (D) / Actual (E) figures for the B440 report, and the
# starting Progress (G5) base value. Formulas in F and G columns (and the
# totals in row 37) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Starting progress baseline (no formula, raw value)
$ws.Range("G5").Value = -2419

# Plan (D) and Actual (E) values per day row
$planActual = @{
    6  = @(110, 117)
    7  = @(110, 114)
    8  = @(110, 106)
    9  = @(110, 119)
    12 = @(110, 125)
    13 = @(144, 120)
    14 = @(110, 110)
    15 = @(110, 104)
    16 = @(110, 114)
    19 = @(110, 122)
    20 = @(144, 148)
    21 = @(110, 114)
    22 = @(110, 135)
    23 = @(110, 87)
    26 = @(110, $null)
    27 = @(110, $null)
    28 = @(110, $null)
    29 = @(110, $null)
    30 = @(110, $null)
    33 = @(110, $null)
    34 = @(110, $null)
    35 = @(41, $null)
}

foreach ($row in $planActual.Keys) {
    $values = $planActual[$row]
    $planVal = $values[0]
    $actualVal = $values[1]

    $ws.Cells.Item($row, 4).Value = $planVal   # column D - Plan
    if ($null -ne $actualVal) {
        $ws.Cells.Item($row, 5).Value = $actualVal  # column E - Actual
    }
}

$wb.Application.Calculate()
